# Optuna Attempt (go back with original)
# Updates the "Forecast Comparison" sheet's Seasonality Index / Inventory
# Coverage / MyForecast figures plus the derived "Summary" totals to match
# the restored (pre-tuning) forecast values.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison" ---

# Row 2 (W8)
$ws1.Range("L2").Value = 0.92

# Row 3 (W9)
$ws1.Range("L3").Value = 0.84

# Row 4 (W10)
$ws1.Range("L4").Value = 1.1

# Row 5 (W11)
$ws1.Range("L5").Value = 1.02

# Row 6 (W12)
$ws1.Range("L6").Value = 1.14

# Row 7 (W13)
$ws1.Range("L7").Value = 1.07

# Row 8 (W14)
$ws1.Range("D8").Value = 1
$ws1.Range("H8").Value = 13.09
$ws1.Range("L8").Value = 0.84

# Row 9 (W15)
$ws1.Range("H9").Value = 10
$ws1.Range("L9").Value = 0.82

# Row 10 (W16)
$ws1.Range("H10").Value = 9
$ws1.Range("L10").Value = 0.86

# Row 11 (W17)
$ws1.Range("H11").Value = 8
$ws1.Range("L11").Value = 0.82

# Row 12 (W18)
$ws1.Range("H12").Value = 7
$ws1.Range("L12").Value = 0.9399999999999999

# Row 13 (W19)
$ws1.Range("H13").Value = 6
$ws1.Range("L13").Value = 0.9399999999999999

# Row 14 (W20)
$ws1.Range("D14").Value = 1
$ws1.Range("H14").Value = 6.05
$ws1.Range("L14").Value = 1.01

# Row 15 (W21)
$ws1.Range("D15").Value = 1
$ws1.Range("H15").Value = 5.05
$ws1.Range("L15").Value = 1.06

# Row 16 (W22)
$ws1.Range("H16").Value = 3.35
$ws1.Range("L16").Value = 0.99

# Row 17 (W23)
$ws1.Range("H17").Value = 2.35
$ws1.Range("L17").Value = 1.03

# --- Sheet "Summary" ---
# B9/B10 are stored as text, not numbers - force Text format first so the
# digit-only strings aren't silently re-typed as numbers on assignment.
$ws2.Range("B9:B10").NumberFormat = "@"
$ws2.Range("B9").Value = "32"
$ws2.Range("B10").Value = "16"
